$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Update the cells that changed value on existing rows 2-28.
# (Column D = Fecha, I = Calidad, J = Volumen, K = Precio minimo,
#  L = Precio maximo, M = Precio promedio ponderado, P = Precio $/Kg)
# ------------------------------------------------------------------
# Row 2
$ws.Range("D2").Value = 44882
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 700
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = 750
$ws.Range("P2").Value = 750

# Row 3
$ws.Range("D3").Value = 44882
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = 600
$ws.Range("P3").Value = 600

# Row 4
$ws.Range("D4").Value = 44608
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("L4").Value = 650
$ws.Range("M4").Value = 625
$ws.Range("P4").Value = 625

# Row 5
$ws.Range("D5").Value = 44839
$ws.Range("J5").Value = 240
$ws.Range("K5").Value = 700
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = 750
$ws.Range("P5").Value = 750

# Row 6
$ws.Range("D6").Value = 44839
$ws.Range("I6").Value = "Segunda"
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = 600
$ws.Range("P6").Value = 600

# Row 7
$ws.Range("D7").Value = 44624
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 650
$ws.Range("L7").Value = 700
$ws.Range("M7").Value = 675
$ws.Range("P7").Value = 675

# Row 8
$ws.Range("D8").Value = 44764
$ws.Range("J8").Value = 200

# Row 9
$ws.Range("D9").Value = 44764
$ws.Range("J9").Value = 150

# Row 10
$ws.Range("D10").Value = 44859
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 700
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = 750
$ws.Range("P10").Value = 750

# Row 11
$ws.Range("D11").Value = 44859
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = 600
$ws.Range("P11").Value = 600

# Row 12
$ws.Range("D12").Value = 44804
$ws.Range("K12").Value = 750
$ws.Range("L12").Value = 850
$ws.Range("M12").Value = 800
$ws.Range("P12").Value = 800

# Row 13
$ws.Range("D13").Value = 44804
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 650
$ws.Range("L13").Value = 650
$ws.Range("M13").Value = 650
$ws.Range("P13").Value = 650

# Row 14
$ws.Range("D14").Value = 44761
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 700
$ws.Range("L14").Value = 800
$ws.Range("M14").Value = 750
$ws.Range("P14").Value = 750

# Row 15
$ws.Range("D15").Value = 44761
$ws.Range("J15").Value = 150
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = 600
$ws.Range("P15").Value = 600

# Row 16
$ws.Range("D16").Value = 44797
$ws.Range("J16").Value = 240

# Row 17
$ws.Range("D17").Value = 44797
$ws.Range("I17").Value = "Segunda"
$ws.Range("K17").Value = 650
$ws.Range("L17").Value = 650
$ws.Range("M17").Value = 650
$ws.Range("P17").Value = 650

# Row 18
$ws.Range("D18").Value = 44811
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 750
$ws.Range("L18").Value = 850
$ws.Range("M18").Value = 800
$ws.Range("P18").Value = 800

# Row 19
$ws.Range("D19").Value = 44868

# Row 20
$ws.Range("D20").Value = 44754
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 700
$ws.Range("L20").Value = 750
$ws.Range("M20").Value = 725
$ws.Range("P20").Value = 725

# Row 21
$ws.Range("D21").Value = 44837
$ws.Range("J21").Value = 200

# Row 22
$ws.Range("D22").Value = 44837
$ws.Range("J22").Value = 150

# Row 23
$ws.Range("D23").Value = 44831
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 700
$ws.Range("M23").Value = 750
$ws.Range("P23").Value = 750

# Row 24
$ws.Range("D24").Value = 44831
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 600
$ws.Range("P24").Value = 600

# Row 25
$ws.Range("D25").Value = 44791
$ws.Range("J25").Value = 240
$ws.Range("K25").Value = 750
$ws.Range("L25").Value = 800
$ws.Range("M25").Value = 775
$ws.Range("P25").Value = 775

# Row 26
$ws.Range("D26").Value = 44791
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 250
$ws.Range("K26").Value = 650
$ws.Range("L26").Value = 650
$ws.Range("M26").Value = 650
$ws.Range("P26").Value = 650

# Row 27
$ws.Range("D27").Value = 44610
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 650
$ws.Range("M27").Value = 625
$ws.Range("P27").Value = 625

# Row 28
$ws.Range("D28").Value = 44818
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 800
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = 850
$ws.Range("P28").Value = 850

# ------------------------------------------------------------------
# Two new data rows (29 and 30) appended at the bottom, extending
# the sheet dimension from A1:R28 to A1:R30.
# ------------------------------------------------------------------
# Row 29
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 44799
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = 100112044
$ws.Range("G29").Value = "Perejil"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 160
$ws.Range("K29").Value = 750
$ws.Range("L29").Value = 850
$ws.Range("M29").Value = 800
$ws.Range("N29").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 800
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
$ws.Range("D29").NumberFormat = $ws.Range("D2").NumberFormat

# Row 30
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = 44799
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = 100112044
$ws.Range("G30").Value = "Perejil"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Segunda"
$ws.Range("J30").Value = 120
$ws.Range("K30").Value = 650
$ws.Range("L30").Value = 650
$ws.Range("M30").Value = 650
$ws.Range("N30").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 650
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"
$ws.Range("D30").NumberFormat = $ws.Range("D2").NumberFormat

Write-Output "Perejil sheet updated: rows 2-28 refreshed, rows 29-30 added."
